$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (Text number format) on specific D-column cells whose new values
# would otherwise be auto-parsed by Excel as numbers (losing the original
# string formatting, e.g. trailing zeros or precision).
$textForceCells = @("D5","D10","D14","D15","D16","D20","D22","D23","D24","D25","D27","D29","D32","D34","D39","D40","D42","D43","D46","D47","D49","D50","D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values row by row, in original top-to-bottom sheet order.
# Row 2
$ws.Range("D2").Value = "26.177.38"
$ws.Range("E2").Value = "  +0.05%  "
# Row 3
$ws.Range("D3").Value = "1.586.06"
$ws.Range("E3").Value = "  +0.62%  "
# Row 4
$ws.Range("E4").Value = "  -0.01%  "
# Row 5
$ws.Range("D5").Value = "211.87"
$ws.Range("E5").Value = "  +1.51%  "
# Row 7
$ws.Range("E7").Value = "  -0.03%  "
# Row 8
$ws.Range("E8").Value = "  +0.42%  "
# Row 9
$ws.Range("E9").Value = "  -0.48%  "
# Row 10
$ws.Range("D10").Value = "19.26"
$ws.Range("E10").Value = "  -1.65%  "
# Row 11
$ws.Range("E11").Value = "  +0.55%  "
# Row 12
$ws.Range("D12").Value = "1.809.80"
# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.620.31"
$ws.Range("E13").Value = "  +2.90%  "
# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "4.01"
$ws.Range("E14").Value = "  -1.21%  "
# Row 15
$ws.Range("D15").Value = "0.519"
$ws.Range("E15").Value = "  +0.74%  "
# Row 16
$ws.Range("D16").Value = "64.06"
$ws.Range("E16").Value = "  -0.52%  "
# Row 17
$ws.Range("D17").Value = "26.184.70"
$ws.Range("E17").Value = "  +0.09%  "
# Row 18
$ws.Range("D18").Value = "0.0₃0725"
$ws.Range("E18").Value = "  -0.31%  "
# Row 19
$ws.Range("E19").Value = "  +1.47%  "
# Row 20
$ws.Range("D20").Value = "212.57"
$ws.Range("E20").Value = "  +1.77%  "
# Row 21
$ws.Range("E21").Value = "  +0.06%  "
# Row 22
$ws.Range("D22").Value = "4.23"
$ws.Range("E22").Value = "  -0.47%  "
# Row 23
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").Value = "2.16"
$ws.Range("E23").Value = "  -0.03%  "
# Row 24
$ws.Range("B24").Value = "Avalanche"
$ws.Range("C24").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D24").Value = "8.98"
$ws.Range("E24").Value = "  +1.67%  "
# Row 25
$ws.Range("D25").Value = "143.53"
$ws.Range("E25").Value = "  -0.23%  "
# Row 26
$ws.Range("E26").Value = "  +0.04%  "
# Row 27
$ws.Range("D27").Value = "6.98"
$ws.Range("E27").Value = "  +0.00%  "
# Row 28
$ws.Range("E28").Value = "  -0.74%  "
# Row 29
$ws.Range("D29").Value = "15.17"
$ws.Range("E29").Value = "  -0.24%  "
# Row 30
$ws.Range("E30").Value = "  -1.93%  "
# Row 31
$ws.Range("E31").Value = "  +1.13%  "
# Row 32
$ws.Range("D32").Value = "3.20"
$ws.Range("E32").Value = "  -0.55%  "
# Row 33
$ws.Range("D33").Value = "1.338.02"
$ws.Range("E33").Value = "  +4.57%  "
# Row 34
$ws.Range("D34").Value = "2.93"
$ws.Range("E34").Value = "  -2.19%  "
# Row 35
$ws.Range("E35").Value = "  +0.14%  "
# Row 36
$ws.Range("E36").Value = "  -0.88%  "
# Row 37
$ws.Range("E37").Value = "  -4.56%  "
# Row 38
$ws.Range("E38").Value = "  +0.29%  "
# Row 39
$ws.Range("D39").Value = "0.818"
$ws.Range("E39").Value = "  +0.99%  "
# Row 40
$ws.Range("D40").Value = "5.79"
$ws.Range("E40").Value = "  +3.61%  "
# Row 41
$ws.Range("E41").Value = "  -0.05%  "
# Row 42
$ws.Range("D42").Value = "0.963"
$ws.Range("E42").Value = "  -12.83%  "
# Row 43
$ws.Range("D43").Value = "0.770"
$ws.Range("E43").Value = "  +0.59%  "
# Row 44
$ws.Range("E44").Value = "  +0.51%  "
# Row 45
$ws.Range("D45").Value = "1.721.82"
$ws.Range("E45").Value = "  +0.59%  "
# Row 46
$ws.Range("D46").Value = "60.99"
$ws.Range("E46").Value = "  -2.23%  "
# Row 47
$ws.Range("D47").Value = "85.83"
$ws.Range("E47").Value = "  -3.24%  "
# Row 48
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0102"
$ws.Range("E48").Value = "  -0.72%  "
# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.48"
$ws.Range("E49").Value = "  -1.52%  "
# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.0980"
$ws.Range("E50").Value = "  -2.63%  "
# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0501"
$ws.Range("E51").Value = "  -0.79%  "
